$d = $word.ActiveDocument

# 5° Mostrar o site  ->  append " – durante isso dizer qual foi a maior dificuldade e superação"
$d.Content.Find.Execute("5° Mostrar o site", $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "5° Mostrar o site – durante isso dizer qual foi a maior dificuldade e superação",
                         2)

# 6° Mostrar o modelo lógico após explicar o sistema de chips
#   -> 6° Mostrar o modelo lógico após mostrar o site
$d.Content.Find.Execute("o modelo lógico após explicar o sistema de chips", $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "o modelo lógico após mostrar o site",
                         2)
